# Adds the "Paginator" reference line (with its hyperlink-styled URL)
# to the empty, underline-formatted paragraph that follows the
# "Validator Contraseña" reference entry, right before the trailing
# blank paragraphs at the end of the document.

$d = $word.ActiveDocument

# --- Locate the target paragraph -------------------------------------
# Anchor on the unique, unambiguous text of the previous reference
# entry ("Validator Contraseña: https://dev.to/.../...-3pkl") and then
# move one paragraph further — that is the empty paragraph (with a
# single-underline run-format left on its mark) that needs the new
# "Paginator" content.
$anchor = $d.Content
$found = $anchor.Find.Execute("custom-validator-3pkl", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text for the Paginator reference paragraph"
}
$anchorParaIndex = $anchor.Paragraphs.Item(1).Index
$targetIndex = $anchorParaIndex + 1
$targetPara = $d.Paragraphs.Item($targetIndex)
$targetRange = $targetPara.Range

# --- Insert the plain-text portion + proofErr spell-check markers ----
# Word flags "Paginator" as a misspelling, wrapping it in
# <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>.
# We rebuild the whole paragraph via InsertXML so those proofErr marks
# come through, while explicitly keeping the paragraph's original
# identity attributes (paraId/textId/rsids) and its <w:pPr> (the
# single-underline run format) untouched.
$paraId = "6A9DBBB7"
$textId = "77777777"
$rsidR = "00472965"

$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="' + $paraId + '" w14:textId="' + $textId + '" w:rsidR="' + $rsidR + '" w:rsidRPr="' + $rsidR + '" w:rsidRDefault="' + $rsidR + '"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Paginator</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetRange.InsertXML($xml)

# --- Append the hyperlink-styled URL run ------------------------------
$targetPara = $d.Paragraphs.Item($targetIndex)
$paraRange = $targetPara.Range
$insertAt = $paraRange.End - 1
$urlText = "https://www.tutorialesprogramacionya.com/angularya/detalleconcepto.php?punto=46&codigo=46&inicio=40"
$ip = $d.Range($insertAt, $insertAt)
$ip.InsertAfter($urlText)

$urlStart = $insertAt
$urlEnd = $urlStart + $urlText.Length
$urlRange = $d.Range($urlStart, $urlEnd)
$urlRange.Style = "Hipervnculo"

Write-Output "Paginator reference inserted"
